# Update Col1a2-Itga2 NATMI output with newly recomputed TPM-based values.
#
# The workbook holds one row per (Sending cluster, Target cluster) pair for
# the Col1a2 -> Itga2 ligand-receptor edge. Per sending cluster there is a
# single "Ligand average expression value" (col G); per target cluster there
# is a single "Receptor average expression value" (col M). Every other
# numeric column in the sheet is mechanically derived from those two base
# numbers, so we only need to supply the refreshed G/M values per cluster and
# recompute the rest the same way NATMI does:
#
#   H (ligand total)   = G * 3                              (3 samples/cluster)
#   I (ligand avg spec) = G / SUM(G over all sending clusters for this edge)
#   J (ligand tot spec) = H / SUM(H over all sending clusters for this edge)
#   N (receptor total)  = M * 3
#   O (receptor avg spec) = M / SUM(M over all target clusters for this edge)
#   P (receptor tot spec) = N / SUM(N over all target clusters for this edge)
#   Q (edge avg weight)  = G * M
#   R (edge tot weight)  = H * N
#   S (edge avg spec)    = I * O
#   T (edge tot spec)    = J * P

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed ligand average expression value per sending cluster (col G).
$ligandAvg = @{
    'ECs'           = 18.31647966666667
    'FAPs'          = 3161.845459
    'MuSCs'         = 155.6514383333333
    'Resolving-Mac' = 2.142642
}

# Refreshed receptor average expression value per target cluster (col M).
$receptorAvg = @{
    'ECs'           = 3.425446666666666
    'FAPs'          = 1.077748
    'MuSCs'         = 0.62317
    'Resolving-Mac' = 0.01852966666666667
}

$sumLigandAvg = 0
foreach ($v in $ligandAvg.Values) { $sumLigandAvg += $v }
$sumLigandTot = $sumLigandAvg * 3

$sumReceptorAvg = 0
foreach ($v in $receptorAvg.Values) { $sumReceptorAvg += $v }
$sumReceptorTot = $sumReceptorAvg * 3

for ($row = 2; $row -le 17; $row++) {
    $sending = $ws.Cells.Item($row, 1).Value2
    $target  = $ws.Cells.Item($row, 4).Value2

    $g = $ligandAvg[$sending]
    $h = $g * 3
    $i = $g / $sumLigandAvg
    $j = $h / $sumLigandTot

    $m = $receptorAvg[$target]
    $n = $m * 3
    $o = $m / $sumReceptorAvg
    $p = $n / $sumReceptorTot

    $q = $g * $m
    $r = $h * $n
    $s = $i * $o
    $t = $j * $p

    $ws.Cells.Item($row, 7).Value  = $g   # G: Ligand average expression value
    $ws.Cells.Item($row, 8).Value  = $h   # H: Ligand total expression value
    $ws.Cells.Item($row, 9).Value  = $i   # I: Ligand derived specificity (avg)
    $ws.Cells.Item($row, 10).Value = $j   # J: Ligand derived specificity (total)

    $ws.Cells.Item($row, 13).Value = $m   # M: Receptor average expression value
    $ws.Cells.Item($row, 14).Value = $n   # N: Receptor total expression value
    $ws.Cells.Item($row, 15).Value = $o   # O: Receptor derived specificity (avg)
    $ws.Cells.Item($row, 16).Value = $p   # P: Receptor derived specificity (total)

    $ws.Cells.Item($row, 17).Value = $q   # Q: Edge average expression weight
    $ws.Cells.Item($row, 18).Value = $r   # R: Edge total expression weight
    $ws.Cells.Item($row, 19).Value = $s   # S: Edge average expression derived specificity
    $ws.Cells.Item($row, 20).Value = $t   # T: Edge total expression derived specificity
}
